$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data (values are stored as text to preserve exact formatting)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.948.31"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.705.64"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.97"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3951"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4028"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.482"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.74"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.57%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08816"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.20"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.470"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.992"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001355"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.716.87"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "96.29"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07185"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.59"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.358"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.942.67"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.24%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.352"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.66"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +5.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.203"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +15.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.77"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "150.44"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +8.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.429"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.627"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +33.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.905.20"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08562"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.20%  "
# Rows 35 and 36: ImmutableX and VeChain swapped ranking positions
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03137"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.48%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.046"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.199"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2859"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09557"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.54%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.00"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.37"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.700"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7393"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.249"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.408"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08792"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +9.18%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "139.20"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.00%  "
